# Update countries & provincias Spain
# Refresh the COVID stats table on sheet "Pais":
#  - header timestamp (A1) bumped from 14:22 to 14:52
#  - several countries received updated totals, which shifts their
#    sort position (the sheet is kept sorted descending by column B,
#    "Casos totales"), pushing neighbouring rows down by one position
#    while keeping their previously-reported figures intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: "Datos actualizados a 25 de Abril de 2020 a las 14:22" -> "...14:52"
$ws.Range("A1").Value = "Datos actualizados a 25 de Abril de 2020 a las 14:52"

# --- Cell-level updates (row, col, newValue). col: A=1 .. H=8
$updates = @(
    # India (row 19) - refreshed totals
    @(19, 2, 24942),
    @(19, 3, 495),
    @(19, 5, 18664),

    # Portugal (row 20) - refreshed totals
    @(20, 2, 23392),
    @(20, 3, 595),
    @(20, 4, 1277),
    @(20, 5, 21235),
    @(20, 6, 186),
    @(20, 7, 26),
    @(20, 8, 880),

    # Suecia (row 24) - refreshed totals
    @(24, 2, 18177),
    @(24, 3, 610),
    @(24, 5, 14980),
    @(24, 7, 40),
    @(24, 8, 2192),

    # Row 25 becomes Arabia Saudita (fresh, higher totals than Israel/Austria)
    @(25, 1, 'Arabia Saudita'),
    @(25, 2, 16299),
    @(25, 3, 1197),
    @(25, 4, 2215),
    @(25, 5, 13948),
    @(25, 6, 93),
    @(25, 7, 9),
    @(25, 8, 136),

    # Row 26 becomes Israel, keeping Israel's previously-reported figures
    @(26, 1, 'Israel'),
    @(26, 3, 90),
    @(26, 4, 6159),
    @(26, 5, 8791),
    @(26, 6, 130),
    @(26, 7, 4),
    @(26, 8, 198),

    # Row 27 becomes Austria, keeping Austria's previously-reported figures
    @(27, 1, 'Austria'),
    @(27, 2, 15148),
    @(27, 3, 77),
    @(27, 4, 12103),
    @(27, 5, 2509),
    @(27, 6, 148),
    @(27, 7, 6),
    @(27, 8, 536),

    # Noruega (row 43) - refreshed totals
    @(43, 5, 7231),
    @(43, 6, 50),
    @(43, 7, 1),
    @(43, 8, 200),

    # Kazajistan (row 62) - refreshed totals
    @(62, 2, 2564),
    @(62, 3, 148),
    @(62, 5, 1910),

    # Republica de Yibuti (row 88) - refreshed totals
    @(88, 2, 1008),
    @(88, 3, 9),
    @(88, 4, 373),
    @(88, 5, 633),

    # Row 110 becomes Sri Lanka (fresh, higher totals than Guatemala/Taiwan)
    @(110, 1, 'Sri Lanka'),
    @(110, 2, 433),
    @(110, 3, 16),
    @(110, 4, 116),
    @(110, 5, 310),
    @(110, 6, 2),
    @(110, 8, 7),

    # Row 111 becomes Guatemala, keeping Guatemala's previously-reported figures
    @(111, 1, 'Guatemala'),
    @(111, 2, 430),
    @(111, 3, 46),
    @(111, 4, 30),
    @(111, 5, 389),
    @(111, 6, 5),
    @(111, 8, 11),

    # Row 112 becomes Taiwan, keeping Taiwan's previously-reported figures
    @(112, 1, 'Taiwan'),
    @(112, 2, 429),
    @(112, 3, 1),
    @(112, 4, 275),
    @(112, 5, 148),
    @(112, 6, 0),
    @(112, 8, 6),

    # Mozambique (row 161) - refreshed totals
    @(161, 2, 70),
    @(161, 3, 5),
    @(161, 5, 58)
)

foreach ($u in $updates) {
    $row = $u[0]
    $col = $u[1]
    $val = $u[2]
    $ws.Cells.Item($row, $col).Value = $val
}
